$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.151.30'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.834.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("E6").Value = '  -2.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07432'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2940'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07744'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.874.80'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.95%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.991'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6700'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.103'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008371'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.104.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.065.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '227.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.159'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.636'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1404'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.511'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.116'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.043'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05341'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("E34").Value = '  +1.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7546'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.139'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.671'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.274.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.733'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9292'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.08912'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.971'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.55'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.957.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5157'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("E49").Value = '  -1.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.50'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05909'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.62%  '
